# Add two new worksheets ("Partial" and "Complete") that each carry a
# subset / full copy of Sheet1's data, then switch the active tab to the
# newly added "Complete" sheet.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")

# --- New sheet: "Partial" (columns A:B only, no header styling) ---
$partial = $wb.Worksheets.Add($null, $sheet1)
$partial.Name = "Partial"

$partial.Range("A1").Value = "String Property"
$partial.Range("B1").Value = "Numeric Property"

$partial.Range("A2").Value = "Hello"
$partial.Range("B2").Value = 1

$partial.Range("A3").Value = "Good Morning"
$partial.Range("B3").Value = 2

$partial.Range("A4").Value = "Good Afternoon"
$partial.Range("B4").Value = 3

$partial.Range("A1:C4").Select()

# --- New sheet: "Complete" (columns A:C, full copy of Sheet1's first three columns) ---
$complete = $wb.Worksheets.Add($null, $partial)
$complete.Name = "Complete"

$complete.Range("A1").Value = "String Property"
$complete.Range("B1").Value = "Numeric Property"
$complete.Range("C1").Value = "Last One"

$complete.Range("A2").Value = "Hello"
$complete.Range("B2").Value = 1
$complete.Range("C2").Value = "Goodbye"

$complete.Range("A3").Value = "Good Morning"
$complete.Range("B3").Value = 2
$complete.Range("C3").Value = "Good Night"

$complete.Range("A4").Value = "Good Afternoon"
$complete.Range("B4").Value = 3
$complete.Range("C4").Value = "Good Evening"

$complete.Range("D8").Select()

$complete.Activate()
